# Changes of 31st March 2022
# Updates the ShipmentTrackNum / PackageTrackNum columns (C, D) for rows 2-22
# of Sheet1 with a fresh batch of tracking numbers, while keeping the cell
# type as text (matching the original "t=s" shared-string cells) rather than
# letting Excel auto-coerce the numeric-looking strings into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used as a relay: writing the new value there as a text
# formula result and then Copy / PasteSpecial-values it into the real
# target keeps the destination cell's existing formatting/style untouched
# (a direct .Value assignment of a numeric-looking string would silently
# re-type the cell as a number and requires a NumberFormat change to avoid
# it, which would touch styles.xml).
$scratch = $ws.Range("Z100")

function Set-TextValue($rangeAddr, $value) {
    $scratch.Formula = '=TEXT(' + $value + ',"0")'
    $scratch.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
}

Set-TextValue "C2" "320018191948"
Set-TextValue "C3" "320018191959"
Set-TextValue "C4" "320018191981"
Set-TextValue "C5" "320018192006"
Set-TextValue "D5" "320018192006"
Set-TextValue "C6" "320018192040"
Set-TextValue "D6" "320018192040"
Set-TextValue "C7" "320018192061"
Set-TextValue "D7" "320018192061"
Set-TextValue "C8" "320018192094"
Set-TextValue "C9" "320018192131"
Set-TextValue "C10" "320018192164"
Set-TextValue "C11" "320018192186"
Set-TextValue "C12" "320018192223"
Set-TextValue "C13" "320018192245"
Set-TextValue "D13" "320018192245"
Set-TextValue "C14" "320018192278"
Set-TextValue "D14" "320018192278"
Set-TextValue "C15" "320018192290"
Set-TextValue "D15" "320018192290"
Set-TextValue "C16" "320018192326"
Set-TextValue "D16" "320018192326"
Set-TextValue "C17" "320018192348"
Set-TextValue "D17" "320018192348"
Set-TextValue "C18" "320018192381"
Set-TextValue "C19" "320018192407"
Set-TextValue "C20" "320018192430"
Set-TextValue "C21" "320018192451"
Set-TextValue "C22" "320018192484"

$scratch.ClearContents()
$excel.CutCopyMode = 0
